$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill empty predicted price: compute TimeTaken in Hours from TimeTaken in Minutes
$ws.Range("C2").Formula = "=B2/60"
